$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = 46047
$ws.Range("B11").Value = 48994
$ws.Range("A11:B11").NumberFormat = $ws.Range("A10:B10").NumberFormat
